# "Generate Report for Handback"
#
# The handback finished (de-de is now in sync with en-US at
# 2016-08-19 21:02:29 / 21:02:35), so the localization-status report is
# regenerated: the Overview/zh-cn/de-de "Status" cells flip from
# "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" column now links straight to the handed-back
# source markdown file, the "Latest Handback File"/"Latest Handback
# DateTime" columns are filled in, and a few columns are widened so the
# new, longer text fits.

$wb = $excel.ActiveWorkbook

$mdName   = "4987c566-b142-4352-bddb-92d8c3dc69ee.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e420eb3e9991c04099401a7b6071cad3f432208/e2e/4987c566-b142-4352-bddb-92d8c3dc69ee.md"
$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (zh-cn / de-de) now read "Handed back"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# Widen the two status columns so the longer text is fully visible.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Status
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668

# Latest Target File (I2) now links to the handed-back source file.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdName)

# Latest Handback File / Latest Handback DateTime
$wsZhCn.Range("J2").Value = "4987c566-b142-4352-bddb-92d8c3dc69ee.c4809855ef7909218a54bfdcf61b514fc1587d33.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-19 21:02:29"

# Latest Target File / Latest Handback File columns are now wide enough
# to show a full file name.
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668

# Latest Target File (I2) now links to the handed-back source file.
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdName)

# Latest Handback File / Latest Handback DateTime
$wsDeDe.Range("J2").Value = "4987c566-b142-4352-bddb-92d8c3dc69ee.c4809855ef7909218a54bfdcf61b514fc1587d33.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-19 21:02:35"

$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Report regenerated for handback."
